# "started refactoring for #3"
#
# Renames the sheet, adds a named range over a small lookup table that is
# written into J16:K18, and moves the active selection to reflect where the
# user ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xl/workbook.xml: <sheet name="test - Copy" .../> -> <sheet name="test" .../>
$ws.Name = "test"

# New lookup table backing the named range, written column-by-column
# (J16:J18 then K16:K18) so values line up with J16/K16 sharing the same
# header text "nam_ran_col1".
$ws.Range("J16").Value = "nam_ran_col1"
$ws.Range("J17").Value = "val1"
$ws.Range("J18").Value = "val2"
$ws.Range("K16").Value = "nam_ran_col1"
$ws.Range("K17").Value = "val3"
$ws.Range("K18").Value = "val4"

# xl/workbook.xml: new <definedNames><definedName name="xlsx_named_range1">
# test!$J$16:$K$18</definedName></definedNames>
$wb.Names.Add("xlsx_named_range1", "=test!`$J`$16:`$K`$18")

# xl/worksheets/sheet1.xml: <selection activeCell="Q13" .../> -> activeCell="X15"
$ws.Range("X15").Select()
